{"js": "// Update the date line and the 25 division problems in the single table.\n// We address every text run by POSITION (paragraph / table cell), never by\n// text search, because several new values equal other old values\n// (e.g. \"45\u00f74=\" -> \"91\u00f77=\" while a separate cell holds \"91\u00f77=\" -> \"97\u00f72=\"),\n// so a naive global find/replace could clobber the wrong occurrence.\n\nconst body = context.document.body;\n\n// 1) Date paragraph (the only paragraph outside the table).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2025-09-22 Monday\", \"Replace\");\n\n// 2) The division-problem table (only one table in the document).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// New values, row-major, for the 5 content rows (5 cells each) of the\n// 5x5 grid of division problems.\nconst newValues = [\n  [\"13\u00f77=\", \"96\u00f74=\", \"91\u00f77=\", \"55\u00f72=\", \"83\u00f78=\"],\n  [\"23\u00f77=\", \"66\u00f74=\", \"30\u00f79=\", \"86\u00f74=\", \"24\u00f79=\"],\n  [\"76\u00f79=\", \"84\u00f76=\", \"79\u00f79=\", \"92\u00f72=\", \"83\u00f77=\"],\n  [\"81\u00f76=\", \"90\u00f75=\", \"66\u00f72=\", \"38\u00f74=\", \"11\u00f72=\"],\n  [\"55\u00f73=\", \"94\u00f73=\", \"20\u00f77=\", \"74\u00f78=\", \"97\u00f72=\"],\n];\n\n// Content rows are rows 0, 4, 8, 12, 16 of the 20-row table (3 blank rows\n// follow each content row).\nconst contentRowIndexes = [0, 4, 8, 12, 16];\n\nfor (let r = 0; r < contentRowIndexes.length; r++) {\n  const rowIdx = contentRowIndexes[r];\n  for (let c = 0; c < newValues[r].length; c++) {\n    table.getCell(rowIdx, c).value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the single table.\n# Every text run is addressed by POSITION (paragraph index / table cell),\n# never by text search-and-replace, because several new values equal other\n# old values (e.g. \"45\u00f74=\" -> \"91\u00f77=\" while a separate cell holds\n# \"91\u00f77=\" -> \"97\u00f72=\"), so a naive global Find/Replace could clobber the\n# wrong occurrence depending on execution order.\n\n$d = $word.ActiveDocument\n\n# 1) Date paragraph (the only paragraph outside the table).\n$d.Paragraphs(1).Range.Text = \"2025-09-22 Monday\"\n\n# 2) The division-problem table (only one table in the document).\n$table = $d.Tables(1)\n\n# New values, row-major, for the 5 content rows (5 cells each) of the\n# 5x5 grid of division problems. Content rows are the table's absolute\n# rows 1, 5, 9, 13, 17 (1-based) -- 3 blank rows follow each content row.\n$newValues = @(\n    @(\"13\u00f77=\", \"96\u00f74=\", \"91\u00f77=\", \"55\u00f72=\", \"83\u00f78=\"),\n    @(\"23\u00f77=\", \"66\u00f74=\", \"30\u00f79=\", \"86\u00f74=\", \"24\u00f79=\"),\n    @(\"76\u00f79=\", \"84\u00f76=\", \"79\u00f79=\", \"92\u00f72=\", \"83\u00f77=\"),\n    @(\"81\u00f76=\", \"90\u00f75=\", \"66\u00f72=\", \"38\u00f74=\", \"11\u00f72=\"),\n    @(\"55\u00f73=\", \"94\u00f73=\", \"20\u00f77=\", \"74\u00f78=\", \"97\u00f72=\")\n)\n\n$contentRows = @(1, 5, 9, 13, 17)\n\nfor ($r = 0; $r -lt $contentRows.Length; $r++) {\n    $rowIdx = $contentRows[$r]\n    for ($c = 1; $c -le 5; $c++) {\n        $table.Cell($rowIdx, $c).Range.Text = $newValues[$r][$c - 1]\n    }\n}\n"}
